$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column I (isSelected), shifting existing
# columns I..M to J..N.
$ws.Range("I1").EntireColumn.Insert()

# New header for the inserted column.
$ws.Range("I1").Value = "mu"

# New "mu" values for rows 2..21.
$muValues = @(
    586893.906377139,
    541758.4117156758,
    508723.5533420356,
    590999.8542688669,
    567429.8798277654,
    555899.6609958746,
    559619.1960827942,
    551980.9383945196,
    524890.893446172,
    539480.4326250888,
    557303.1279787013,
    576913.2581635589,
    550529.3242182946,
    585079.3153675391,
    523035.4235170254,
    555902.1090949007,
    584033.6451052462,
    564384.9960362483,
    556481.6741497983,
    569801.4788333059
)

for ($i = 0; $i -lt $muValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $muValues[$i]
}
